$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New candidate rows appended by the administrator functions
$ws.Range("A8").Value = 6666
$ws.Range("B8").Value = "EL54321"
$ws.Range("C8").Value = "Fung"
$ws.Range("D8").Value = "Ming Kong"
$ws.Range("E8").Value = 90156789

$ws.Range("A9").Value = 6666
$ws.Range("B9").Value = "EL666666"
$ws.Range("C9").Value = "Fung"
$ws.Range("D9").Value = "Steve"
$ws.Range("E9").Value = 24484568

$ws.Range("A10").Value = 5555
$ws.Range("B10").Value = "EL12345"
$ws.Range("C10").Value = "Chan"
$ws.Range("D10").Value = "Tai Man"
$ws.Range("E10").Value = 98765432

$ws.Range("A11").Value = 5555
$ws.Range("B11").Value = "EL98765"
$ws.Range("C11").Value = "Wong"
$ws.Range("D11").Value = "Tai Sin"
$ws.Range("E11").Value = 65432109

$ws.Range("E8:E11").NumberFormat = "@"

[void]$ws.Range("B10:E11").Select()
